$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2241.6667
$ws.Range("I40").Value = 2437.625
$ws.Range("J40").Value = 1849.75
$ws.Range("K40").Value = 2437.625
$ws.Range("L40").Value = 1849.75
$ws.Range("M40").Value = -2262.625
$ws.Range("N40").Value = -2199.75
# Row 64
$ws.Range("H64").Value = 3654.7222
$ws.Range("I64").Value = 3818.5
$ws.Range("J64").Value = 3450
$ws.Range("K64").Value = 3818.5
$ws.Range("L64").Value = 3450
$ws.Range("M64").Value = -3570.5
$ws.Range("N64").Value = -3946
# Row 67
$ws.Range("H67").Value = 3654.7222
$ws.Range("I67").Value = 3818.5
$ws.Range("J67").Value = 3450
$ws.Range("K67").Value = 3818.5
$ws.Range("L67").Value = 3450
$ws.Range("M67").Value = -2960.5
$ws.Range("N67").Value = -5166
# Row 76
$ws.Range("H76").Value = 3047.963
$ws.Range("I76").Value = 2695.2273
$ws.Range("K76").Value = 2695.2273
$ws.Range("M76").Value = -2380.2273
# Row 79
$ws.Range("H79").Value = 3047.963
$ws.Range("I79").Value = 2695.2273
$ws.Range("K79").Value = 2695.2273
$ws.Range("M79").Value = -1603.2273
# Row 86
$ws.Range("H86").Value = 3948.6667
$ws.Range("I86").Value = 6458.8
$ws.Range("J86").Value = 2983.2307
$ws.Range("K86").Value = 6458.8
$ws.Range("L86").Value = 2983.2307
$ws.Range("M86").Value = -5335.8
$ws.Range("N86").Value = -5229.2307
# Row 89
$ws.Range("H89").Value = 3948.6667
$ws.Range("I89").Value = 6458.8
$ws.Range("J89").Value = 2983.2307
$ws.Range("K89").Value = 32294
$ws.Range("L89").Value = 14916.1535
$ws.Range("M89").Value = -26678
$ws.Range("N89").Value = -26148.1535
# Row 106
$ws.Range("H106").Value = 4275
$ws.Range("I106").Value = 3509.0908
$ws.Range("J106").Value = 5960
$ws.Range("K106").Value = 3509.0908
$ws.Range("L106").Value = 5960
$ws.Range("M106").Value = -2878.0908
$ws.Range("N106").Value = -7222
# Row 113
$ws.Range("H113").Value = 858308.5600000001
$ws.Range("I113").Value = 3706037
$ws.Range("K113").Value = 3706037
$ws.Range("M113").Value = -3702783
# Row 135
$ws.Range("H135").Value = 740.82355
$ws.Range("I135").Value = 706.26666
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 6356.39994
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -3821.39994
$ws.Range("N135").Value = -14070
# Row 137
$ws.Range("H137").Value = 1809.6666
$ws.Range("I137").Value = 996.4167
$ws.Range("J137").Value = 2351.8333
$ws.Range("K137").Value = 2989.2501
$ws.Range("L137").Value = 7055.499899999999
$ws.Range("M137").Value = -439.2501000000002
$ws.Range("N137").Value = -12155.4999
# Row 138
$ws.Range("H138").Value = 3244.3906
$ws.Range("I138").Value = 1371
$ws.Range("J138").Value = 3676.7114
$ws.Range("K138").Value = 4113
$ws.Range("L138").Value = 11030.1342
$ws.Range("M138").Value = 1027
$ws.Range("N138").Value = -21310.1342

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2165.0557
$ws.Range("I61").Value = 1728.5385
$ws.Range("K61").Value = 1728.5385
$ws.Range("M61").Value = -1516.5385
# Row 63
$ws.Range("H63").Value = 2471.1428
$ws.Range("J63").Value = 2574.5
$ws.Range("L63").Value = 2574.5
$ws.Range("N63").Value = -3946.5
# Row 66
$ws.Range("H66").Value = 2471.1428
$ws.Range("J66").Value = 2574.5
$ws.Range("L66").Value = 12872.5
$ws.Range("N66").Value = -19736.5
# Row 136
$ws.Range("H136").Value = 2165.0557
$ws.Range("I136").Value = 1728.5385
$ws.Range("K136").Value = 5185.6155
$ws.Range("M136").Value = -2635.6155
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# Row 138
$ws.Range("H138").Value = 54429
$ws.Range("J138").Value = 54429
$ws.Range("L138").Value = 54429
$ws.Range("N138").Value = -64709
# Row 139
$ws.Range("H139").Value = 37100
$ws.Range("J139").Value = 37100
$ws.Range("L139").Value = 37100
$ws.Range("N139").Value = -47380
# Row 140
$ws.Range("H140").Value = 74662.5
$ws.Range("J140").Value = 74662.5
$ws.Range("L140").Value = 74662.5
$ws.Range("N140").Value = -85022.5
# Row 141
$ws.Range("H141").Value = 70476.336
$ws.Range("J141").Value = 70476.336
$ws.Range("L141").Value = 70476.336
$ws.Range("N141").Value = -80836.336

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2275033
$ws.Range("I105").Value = 3789955.2
$ws.Range("J105").Value = 2650
$ws.Range("K105").Value = 3789955.2
$ws.Range("L105").Value = 2650
$ws.Range("M105").Value = -3788208.2
$ws.Range("N105").Value = -6144

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5043.263
$ws.Range("I70").Value = 4787
$ws.Range("J70").Value = 6004.25
$ws.Range("K70").Value = 4787
$ws.Range("L70").Value = 6004.25
$ws.Range("M70").Value = -4517
$ws.Range("N70").Value = -6544.25
# Row 73
$ws.Range("H73").Value = 5043.263
$ws.Range("I73").Value = 4787
$ws.Range("J73").Value = 6004.25
$ws.Range("K73").Value = 4787
$ws.Range("L73").Value = 6004.25
$ws.Range("M73").Value = -3851
$ws.Range("N73").Value = -7876.25
# Row 80
$ws.Range("H80").Value = 75256.78999999999
$ws.Range("I80").Value = 2748.1667
$ws.Range("J80").Value = 129638.25
$ws.Range("K80").Value = 2748.1667
$ws.Range("L80").Value = 129638.25
$ws.Range("M80").Value = -1750.1667
$ws.Range("N80").Value = -131634.25
# Row 83
$ws.Range("H83").Value = 75256.78999999999
$ws.Range("I83").Value = 2748.1667
$ws.Range("J83").Value = 129638.25
$ws.Range("K83").Value = 13740.8335
$ws.Range("L83").Value = 648191.25
$ws.Range("M83").Value = -8748.833500000001
$ws.Range("N83").Value = -658175.25
# Row 137
$ws.Range("H137").Value = 35697.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 35697.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 35697.5
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -45897.5
# Row 139
$ws.Range("H139").Value = 44571.285
$ws.Range("J139").Value = 44571.285
$ws.Range("L139").Value = 44571.285
$ws.Range("N139").Value = -54851.285
# Row 140
$ws.Range("H140").Value = 44520
$ws.Range("J140").Value = 43230
$ws.Range("L140").Value = 43230
$ws.Range("N140").Value = -53590
# Row 141
$ws.Range("H141").Value = 63331.5
$ws.Range("J141").Value = 63331.5
$ws.Range("L141").Value = 63331.5
$ws.Range("N141").Value = -73691.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2713.3572
$ws.Range("I40").Value = 1994
$ws.Range("J40").Value = 3672.5
$ws.Range("K40").Value = 1994
$ws.Range("L40").Value = 3672.5
$ws.Range("M40").Value = -1858
$ws.Range("N40").Value = -3944.5
